$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that previously held duplicate/benchmark pricing text
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("J13").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("J14").ClearContents()

# Populate previously empty cells with newly reported benchmark values
$ws.Range("F24").Value = "457,14 TL"
$ws.Range("F25").Value = "380,95 TL"
